$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit numeric-looking string. Assigning it directly through
# .Value would make Excel auto-convert it to a number (General format,
# same as typing it into the cell). To keep it as text (matching the
# original inline-string cell type) without touching the cell's style,
# write it as a formula that evaluates to a text string, then convert the
# formula to a static value via Copy / PasteSpecial(Values only).
$ws.Range("B3").Formula = "=""2570314725427075"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 26.06.2025"

$ws.Range("B6").Value = "27.06."
$ws.Range("C6").Value = "28.06."
$ws.Range("D6").Value = "PAYPAL VKCMOY"
$ws.Range("E6").Value = "88,19-"

$ws.Range("B7").Value = "30.06."
$ws.Range("C7").Value = "01.07."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 69505375"
$ws.Range("E7").Value = "86,11-"

$ws.Range("B8").Value = "04.07."
$ws.Range("C8").Value = "05.07."
$ws.Range("D8").Value = "PAYPAL QVHQCM"
$ws.Range("E8").Value = "49,99-"

$ws.Range("B9").Value = "06.07."
$ws.Range("C9").Value = "07.07."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 70455379"
$ws.Range("E9").Value = "38,86-"

$ws.Range("D12").Value = "KONTOSTAND AM 11.07.2025"
$ws.Range("E12").Value = "263,15-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 20.07.2025"
